# Adds a "common order of use" numbering column (A) to the Scripts sheet,
# indicating the order in which the scripts are meant to be run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scripts")

# Order values per row (row -> order number)
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 5
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A11").Value = 8

# Update the selected cell to reflect where the editor ended up.
$ws.Activate()
$ws.Range("B19").Select()
